$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: new data row for a medicine entry
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "ERASTAPEX TRIO 5/20/12.5MG 30 F.C. TAB"
$ws.Range("H7").Value = "0:2"
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = "1"
$ws.Range("N7").NumberFormat = "@"
$ws.Range("N7").Value = "114.00"
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "228.0000"
$ws.Range("Q7").Value = "2:0"

# Row 8: total row
$ws.Range("N8").Value = 228
$ws.Rows.Item(8).RowHeight = 25.5
